{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Goal (per the diff):\n//  - After the paragraph that reads exactly \"gradlew desktop:dist\"\n//    (the \"Build artifact:\" instructions), insert a brand-new list\n//    paragraph with the text \"The artifact is \" that:\n//      * uses the \"List Paragraph\" style (pStyle -> Listenabsatz)\n//      * is a bulleted list item (numPr ilvl=0 / a fresh numId)\n//      * ends with the (moved) \"_GoBack\" bookmark right after the text\n//  - The original paragraph keeps its proofErr spellcheck wrapper but\n//    loses the bookmark (which now lives in the new paragraph).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Find the exact paragraph \"gradlew desktop:dist\" (not the later\n// \"./gradlew desktop:dist\" code-block sample that appears further\n// down in the document).\nlet sourcePara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === \"gradlew desktop:dist\") {\n    sourcePara = paragraphs.items[i];\n    break;\n  }\n}\nif (!sourcePara) {\n  throw new Error('Could not find the \"gradlew desktop:dist\" paragraph.');\n}\n\n// The existing \"_GoBack\" bookmark sits at the end of this paragraph;\n// it needs to move to the end of the newly typed text, so drop it here.\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// Insert the new paragraph right after it. A trailing sentinel\n// character (\"X\") is typed along with the text so we can anchor the\n// bookmark precisely at \"end of visible text\" before removing the\n// sentinel - inserting the bookmark directly at the paragraph's bare\n// end lands it across the paragraph boundary instead.\nconst newPara = sourcePara.insertParagraph(\"The artifact is X\", Word.InsertLocation.after);\nnewPara.style = \"List Paragraph\";\n\nconst sentinelResults = newPara.search(\"X\", { matchCase: true });\nsentinelResults.load(\"items\");\nawait context.sync();\n\nconst sentinelRange = sentinelResults.items[0];\nconst bookmarkSpot = sentinelRange.getRange(\"Start\");\nbookmarkSpot.insertBookmark(\"_GoBack\");\nsentinelRange.insertText(\"\", Word.InsertLocation.replace);\n\n// Turn the paragraph into a bulleted list item.\nnewPara.startNewList();\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $d (= $word.ActiveDocument) is the open document.\n#\n# Goal (per the diff):\n#  - After the paragraph that reads exactly \"gradlew desktop:dist\"\n#    (the \"Build artifact:\" instructions), insert a brand-new list\n#    paragraph with the text \"The artifact is \" that:\n#      * uses the \"List Paragraph\" style (pStyle -> Listenabsatz)\n#      * is a bulleted list item (ListFormat.ApplyBulletDefault)\n#      * ends with the (moved) \"_GoBack\" bookmark right after the text\n#  - The original paragraph keeps its proofErr spellcheck wrapper but\n#    loses the bookmark (which now lives in the new paragraph).\n\n$d = $word.ActiveDocument\n\n# Find the exact paragraph \"gradlew desktop:dist\" (not the later\n# \"./gradlew desktop:dist\" code-block sample that appears further\n# down in the document).\n$sourcePara = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.Trim()\n    if ($t -eq \"gradlew desktop:dist\") {\n        $sourcePara = $p\n        break\n    }\n}\nif ($sourcePara -eq $null) {\n    throw \"Could not find the 'gradlew desktop:dist' paragraph.\"\n}\n\n# The existing \"_GoBack\" bookmark sits at the end of this paragraph;\n# remove it here since it moves to the end of the newly typed text.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# Insert the new paragraph right after it.\n$srcRange = $sourcePara.Range\n$srcRange.Collapse(0)           # wdCollapseEnd\n$srcRange.InsertParagraphAfter()\n$newPara = $sourcePara.Next()\n\n# Type the text with a trailing sentinel character (\"X\") so we can\n# anchor the bookmark precisely at \"end of visible text\" before\n# removing the sentinel - adding the bookmark directly at the bare\n# paragraph end lands it across the paragraph boundary instead.\n$newPara.Range.Text = \"The artifact is X\"\n\n# Apply the list style/formatting.\n$newPara.Range.Style = \"List Paragraph\"\n$newPara.Range.ListFormat.ApplyBulletDefault()\n\n# Re-fetch the (now list-formatted) paragraph range and locate the\n# sentinel to park the bookmark right before it, then drop the\n# sentinel character.\n$fullRange = $newPara.Range\n$sentinelPos = $fullRange.End - 2\n$bookmarkSpot = $d.Range($sentinelPos, $sentinelPos)\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkSpot)\n$sentinelRange = $d.Range($sentinelPos, $sentinelPos + 1)\n$sentinelRange.Delete()\n"}
